# Generate Report for Handoff
#
# A fresh handoff run replaced the source file's token (GUID-derived name,
# was 411bb499-1b4a-4b5b-b81f-3ac00719a47f) and regenerated the translation
# package content hash (was 4b5e592f3834a347ef5297c1017f69b26cf3a851), so
# every place in the report that echoes the handoff file name / target file
# name / handoff timestamp needs to move forward to the new run's values.

$wb = $excel.ActiveWorkbook

$newGuid = "91123cf2-79a2-4c39-83b8-7388e0c98d98"
$newHash = "7eec8f5ff5d13358a5567e8ec8e862bdf834a4b3"

$newMdName    = "$newGuid.md"
$newZhCnName  = "$newGuid.$newHash.zh-cn.xlf"
$newDeDeName  = "$newGuid.$newHash.de-de.xlf"

# Updates a cell's value and, if that cell carries a hyperlink, keeps the
# hyperlink's displayed text in sync with it (looping over the Hyperlinks
# collection rather than indexing a single Item(n), which - for this
# runtime - is what lets the existing <hyperlink> entry be edited in place
# instead of a stray duplicate getting appended).
function Update-CellAndHyperlink($ws, $cellRef, $newValue) {
    $ws.Range($cellRef).Value = $newValue
    $target = $ws.Range($cellRef).Address()
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $target) {
            $h.TextToDisplay = $newValue
        }
    }
}

# --- Overview sheet: A2 = handoff file, D2 = latest handoff date ---
$wsOverview = $wb.Worksheets.Item("Overview")
Update-CellAndHyperlink $wsOverview "A2" $newMdName
$wsOverview.Range("D2").Value = "2016-03-21 10:54:41"

# --- zh-cn sheet: A2 = handoff file, D2 = target xlf, E2 = handoff datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-CellAndHyperlink $wsZhCn "A2" $newMdName
Update-CellAndHyperlink $wsZhCn "D2" $newZhCnName
$wsZhCn.Range("E2").Value = "2016-03-21 10:54:37"

# --- de-de sheet: A2 = handoff file, D2 = target xlf, E2 = handoff datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
Update-CellAndHyperlink $wsDeDe "A2" $newMdName
Update-CellAndHyperlink $wsDeDe "D2" $newDeDeName
$wsDeDe.Range("E2").Value = "2016-03-21 10:54:41"
